# Daily attendance processing - 2026-02-01 16:45:09
# Swap the "Recorded By" value from "System, dnasr281@gmail.com"
# to "dnasr281@gmail.com, System" wherever it appears (column G,
# "Recorded By"), leaving every other cell (including other G-column
# variants such as just "System" or just "dnasr281@gmail.com") untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$changed = 0

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
        $changed = $changed + 1
    }
}

Write-Host "Updated $changed 'Recorded By' cell(s) in column G."
